$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comment text added in H4 (new shared string)
$ws.Range("H4").Value = "(Tak się liczby dobrały, że jakoś udało mu się skorygować, ale źle)"

# Row 5: D5 becomes 1000 with the same green "OK" fill used elsewhere (e.g. D3)
$ws.Range("D5").Value = 1000
$ws.Range("D5").Interior.Color = $ws.Range("D3").Interior.Color

# Row 5: F5 is cleared out (value + its previous orange fill removed)
$ws.Range("F5").ClearContents()
$ws.Range("F5").ClearFormats()

# Update the active selection to D7
$ws.Range("D7").Select()
